$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New monthly rows to append, matching the existing table layout:
# Serie, Total, 3 meses, 6 meses, 9 meses, 12 meses, 18 meses, 2 años, 5 años, 10 años y más
$newRows = @(
    @("01-07-2021", 6550, 362, 1716, 319, 465, 402, 1406, 1074, 806),
    @("01-08-2021", 8386, 821, 1584, 804, 1865, 732, 1740, 501, 339)
)

$startRow = $ws.Cells.Item($ws.UsedRange.Rows.Count, 1).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A holds a "dd-mm-yyyy" label that must stay plain text (like
    # the existing Serie values), not get auto-converted into a date
    # serial number. Writing it with a leading apostrophe forces Excel to
    # store it as text; resetting the cell style back to Normal afterwards
    # drops the transient quote-prefix formatting so the cell keeps the
    # plain, unstyled look used by every other row in the table.
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Formula = "'" + $rowData[0]
    $aCell.Style = "Normal"

    for ($c = 1; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
